$wb = $excel.ActiveWorkbook

# ALC!row19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 460.33334
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 460.33334
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 460.33334
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -810.33334

# ALC!row41
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 419.64285
$ws.Range("I41").Value = 90
$ws.Range("K41").Value = 90
$ws.Range("M41").Value = 350

# ALC!row43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 922.2
$ws.Range("J43").Value = 922.2
$ws.Range("L43").Value = 922.2
$ws.Range("N43").Value = -1060.2

# ALC!row64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3854.5454
$ws.Range("J64").Value = 4628.5713
$ws.Range("L64").Value = 4628.5713
$ws.Range("N64").Value = -5124.5713

# ALC!row67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3854.5454
$ws.Range("J67").Value = 4628.5713
$ws.Range("L67").Value = 4628.5713
$ws.Range("N67").Value = -6344.5713

# ALC!row94
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 34502
$ws.Range("I94").Value = 1750
$ws.Range("J94").Value = 100006
$ws.Range("K94").Value = 1750
$ws.Range("L94").Value = 100006
$ws.Range("M94").Value = -1299
$ws.Range("N94").Value = -100908

# ALC!row129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 746.85
$ws.Range("I129").Value = 360
$ws.Range("J129").Value = 802.1142599999999
$ws.Range("K129").Value = 1080
$ws.Range("L129").Value = 2406.34278
$ws.Range("M129").Value = 3920
$ws.Range("N129").Value = -12406.34278

# ALC!row132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2748.8647
$ws.Range("I132").Value = 3067.516
$ws.Range("J132").Value = 1102.5
$ws.Range("K132").Value = 9202.548000000001
$ws.Range("L132").Value = 3307.5
$ws.Range("M132").Value = -6672.548000000001
$ws.Range("N132").Value = -8367.5

# ALC!row141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3083.7693
$ws.Range("I141").Value = 3496.8
$ws.Range("J141").Value = 2825.625
$ws.Range("K141").Value = 10490.4
$ws.Range("L141").Value = 8476.875
$ws.Range("M141").Value = -5310.400000000001
$ws.Range("N141").Value = -18836.875

# ARM!row52
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

# ARM!row63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2842421
$ws.Range("I63").Value = 1736.7778
$ws.Range("K63").Value = 1736.7778
$ws.Range("M63").Value = -1050.7778

# ARM!row66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2842421
$ws.Range("I66").Value = 1736.7778
$ws.Range("K66").Value = 8683.889000000001
$ws.Range("M66").Value = -5251.889000000001

# ARM!row74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 20835442
$ws.Range("I74").Value = 28573118
$ws.Range("K74").Value = 28573118
$ws.Range("M74").Value = -28572244

# ARM!row77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 20835442
$ws.Range("I77").Value = 28573118
$ws.Range("K77").Value = 142865590
$ws.Range("M77").Value = -142861222

# ARM!row97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I97").Value = 843.41174
$ws.Range("K97").Value = 843.41174
$ws.Range("M97").Value = -347.41174

# ARM!row102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1361.5883
$ws.Range("I102").Value = 1266.8572
$ws.Range("K102").Value = 1266.8572
$ws.Range("M102").Value = 355.1428000000001

# BSM!row94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 573.4838999999999
$ws.Range("I94").Value = 428.6087
$ws.Range("K94").Value = 428.6087
$ws.Range("M94").Value = 22.3913

# BSM!row99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1999.7142
$ws.Range("I99").Value = 2066.5
$ws.Range("J99").Value = 1599
$ws.Range("K99").Value = 2066.5
$ws.Range("L99").Value = 1599
$ws.Range("M99").Value = -568.5
$ws.Range("N99").Value = -4595

# CRP!row99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3056.3103
$ws.Range("I99").Value = 2606.9
$ws.Range("J99").Value = 4055
$ws.Range("K99").Value = 2606.9
$ws.Range("L99").Value = 4055
$ws.Range("M99").Value = -1108.9
$ws.Range("N99").Value = -7051

# CRP!row126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3056.3103
$ws.Range("I126").Value = 2606.9
$ws.Range("J126").Value = 4055
$ws.Range("K126").Value = 7820.700000000001
$ws.Range("L126").Value = 12165
$ws.Range("M126").Value = -5350.700000000001
$ws.Range("N126").Value = -17105

# CUL!row114
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 2534.5
$ws.Range("I114").Value = 2839.4285
$ws.Range("J114").Value = 400
$ws.Range("K114").Value = 8518.2855
$ws.Range("L114").Value = 1200
$ws.Range("M114").Value = -5264.2855
$ws.Range("N114").Value = -7708

# CUL!row131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 723.28
$ws.Range("J131").Value = 727.5051
$ws.Range("L131").Value = 2182.5153
$ws.Range("N131").Value = -12262.5153

# GSM!row102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1829.6364
$ws.Range("I102").Value = 1458.5555
$ws.Range("J102").Value = 3499.5
$ws.Range("K102").Value = 1458.5555
$ws.Range("L102").Value = 3499.5
$ws.Range("M102").Value = 163.4445000000001
$ws.Range("N102").Value = -6743.5

# GSM!row105
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

# GSM!row139
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H139").Value = 26119.215
$ws.Range("J139").Value = 26119.215
$ws.Range("L139").Value = 26119.215
$ws.Range("N139").Value = -36399.215

# LTW!row100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1961.7646
$ws.Range("I100").Value = 1328.1111
$ws.Range("K100").Value = 1328.1111
$ws.Range("M100").Value = -787.1111000000001

# LTW!row132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1763.1
$ws.Range("I132").Value = 1129.7916
$ws.Range("K132").Value = 3389.3748
$ws.Range("M132").Value = -859.3748000000001

# WVR!row100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 207.92857
$ws.Range("I100").Value = 182.63637
$ws.Range("K100").Value = 365.27274
$ws.Range("M100").Value = 175.72726

# WVR!row122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1601.7368
$ws.Range("I122").Value = 1541.9333
$ws.Range("J122").Value = 1826
$ws.Range("K122").Value = 4625.7999
$ws.Range("L122").Value = 5478
$ws.Range("M122").Value = -2175.7999
$ws.Range("N122").Value = -10378

# WVR!row132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1072.2954
$ws.Range("I132").Value = 868.0741
$ws.Range("J132").Value = 1396.6471
$ws.Range("K132").Value = 2604.2223
$ws.Range("L132").Value = 4189.9413
$ws.Range("M132").Value = -74.22230000000036
$ws.Range("N132").Value = -9249.941299999999

# WVR!row136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 16951084
$ws.Range("I136").Value = 20408944
$ws.Range("J136").Value = 7572
$ws.Range("K136").Value = 61226832
$ws.Range("L136").Value = 22716
$ws.Range("M136").Value = -61224282
$ws.Range("N136").Value = -27816
